$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.649422333333334
$ws.Range("H2").Value = 22.948267
$ws.Range("I2").Value = 0.004484559810904267
$ws.Range("J2").Value = 0.004484559810904268
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 3.740990789035778
$ws.Range("R2").Value = 33.668917101322
$ws.Range("S2").Value = 0.004325123207576119
$ws.Range("T2").Value = 0.004325123207576119

# Row 3
$ws.Range("G3").Value = 7.649422333333334
$ws.Range("H3").Value = 22.948267
$ws.Range("I3").Value = 0.004484559810904267
$ws.Range("J3").Value = 0.004484559810904268
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.1379037858253333
$ws.Range("R3").Value = 1.241134072428
$ws.Range("S3").Value = 0.0001594366033281488
$ws.Range("T3").Value = 0.0001594366033281489

# Row 4
$ws.Range("I4").Value = 0.8893308176045429
$ws.Range("J4").Value = 0.889330817604543
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 741.8740160348982
$ws.Range("R4").Value = 6676.866144314084
$ws.Range("S4").Value = 0.8577130243823085
$ws.Range("T4").Value = 0.8577130243823086

# Row 5
$ws.Range("I5").Value = 0.8893308176045429
$ws.Range("J5").Value = 0.889330817604543
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.03161779322223441
$ws.Range("T5").Value = 0.03161779322223442

# Row 6
$ws.Range("I6").Value = 0.1061846225845528
$ws.Range("J6").Value = 0.1061846225845528
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 88.57852537949621
$ws.Range("R6").Value = 797.2067284154659
$ws.Range("S6").Value = 0.1024095105859576
$ws.Range("T6").Value = 0.1024095105859576

# Row 7
$ws.Range("I7").Value = 0.1061846225845528
$ws.Range("J7").Value = 0.1061846225845528
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("S7").Value = 0.003775111998595203
$ws.Range("T7").Value = 0.003775111998595204
